$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'46.704.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.22%  "
$ws.Range("D3").Value = "'2.311.01"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.28%  "
$ws.Range("E4").Value = "  -0.74%  "
$ws.Range("D5").Value = "'301.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("D6").Value = "'102.81"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +15.62%  "
$ws.Range("D7").Value = "'0.573"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.35%  "
$ws.Range("E8").Value = "  -0.56%  "
$ws.Range("D9").Value = "'0.531"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.93%  "
$ws.Range("D10").Value = "'37.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.49%  "
$ws.Range("D11").Value = "'0.0807"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.61%  "
$ws.Range("D12").Value = "'7.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +8.81%  "
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").Value = "'2.664.44"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.22%  "
$ws.Range("D15").Value = "'2.310.93"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "'14.08"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.35%  "
$ws.Range("D17").Value = "'0.824"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.54%  "
$ws.Range("D18").Value = "'46.693.77"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.90%  "
$ws.Range("D19").Value = "'13.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +22.76%  "
$ws.Range("D20").Value = "'0.0₃0949"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.86%  "
$ws.Range("D21").Value = "'6.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.45%  "
$ws.Range("D22").Value = "'67.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +6.22%  "
$ws.Range("D23").Value = "'248.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.64%  "
$ws.Range("D24").Value = "'2.95"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.06%  "
$ws.Range("D25").Value = "'1.97"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.84%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "
$ws.Range("D27").Value = "'43.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +20.82%  "
$ws.Range("E28").Value = "  +1.69%  "
$ws.Range("D29").Value = "'10.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.28%  "
$ws.Range("D30").Value = "'20.20"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.25%  "
$ws.Range("D31").Value = "'5.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +10.15%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").Value = "'0.0806"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.05%  "
$ws.Range("B33").Value = "Monero"
$ws.Range("C33").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D33").Value = "'146.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.31%  "
$ws.Range("E34").Value = "  +4.58%  "
$ws.Range("D35").Value = "'3.13"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.24%  "
$ws.Range("D36").Value = "'0.112"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.34%  "
$ws.Range("E37").Value = "  +3.65%  "
$ws.Range("E38").Value = "  +9.99%  "
$ws.Range("D39").Value = "'15.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +21.42%  "
$ws.Range("D40").Value = "'4.14"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +18.05%  "
$ws.Range("D41").Value = "'3.46"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +13.01%  "
$ws.Range("D42").Value = "'0.0306"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.33%  "
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").Value = "'1.98"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +18.46%  "
$ws.Range("D45").Value = "'1.854.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +3.60%  "
$ws.Range("D46").Value = "'89.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +23.02%  "
$ws.Range("D47").Value = "'0.196"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.94%  "
$ws.Range("D48").Value = "'74.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +14.59%  "
$ws.Range("D49").Value = "'4.95"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +10.94%  "
$ws.Range("D50").Value = "'97.58"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.29%  "
$ws.Range("D51").Value = "'54.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +11.00%  "
